# Risultati.xlsx - "Relazione 95% e risultati" edit
#
# The four result tables on Foglio1 each had a row-label column (A, then G
# for the second sub-table) plus a duplicate "F1" header cell in column B/H.
# The edit drops the "F1" column entirely and moves the N16/N32/N64/N128 /
# G16/G32/G64/G128 row labels from column A/G into column B/H (replacing
# the numeric values that used to sit there), then removes the now-empty
# A and G columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the B/H numeric cells with the row-label text that used to
#     live in column A / G (for all four tables). ---
$ws.Range("B3").Value  = "N16"
$ws.Range("H3").Value  = "N16"
$ws.Range("B4").Value  = "N32"
$ws.Range("H4").Value  = "N32"
$ws.Range("B5").Value  = "N64"
$ws.Range("H5").Value  = "N64"
$ws.Range("B6").Value  = "N128"
$ws.Range("H6").Value  = "N128"

$ws.Range("B12").Value = "G16"
$ws.Range("H12").Value = "G16"
$ws.Range("B13").Value = "G32"
$ws.Range("H13").Value = "G32"
$ws.Range("B14").Value = "G64"
$ws.Range("H14").Value = "G64"
$ws.Range("B15").Value = "G128"
$ws.Range("H15").Value = "G128"

$ws.Range("B25").Value = "N16"
$ws.Range("H25").Value = "N16"
$ws.Range("B26").Value = "N32"
$ws.Range("H26").Value = "N32"
$ws.Range("B27").Value = "N64"
$ws.Range("H27").Value = "N64"
$ws.Range("B28").Value = "N128"
$ws.Range("H28").Value = "N128"

$ws.Range("B34").Value = "G16"
$ws.Range("H34").Value = "G16"
$ws.Range("B35").Value = "G32"
$ws.Range("H35").Value = "G32"
$ws.Range("B36").Value = "G64"
$ws.Range("H36").Value = "G64"
$ws.Range("B37").Value = "G128"
$ws.Range("H37").Value = "G128"

# --- Drop the now-redundant "F1" header cells (B2/H2 and the matching
#     rows for the other three tables). ---
$ws.Range("B2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("H11").ClearContents()
$ws.Range("B24").ClearContents()
$ws.Range("H24").ClearContents()
$ws.Range("B33").ClearContents()
$ws.Range("H33").ClearContents()

# --- Column A and column G only ever held the row labels that have now
#     been folded into B/H, so clear them completely. ---
$ws.Range("A1:A37").ClearContents()
$ws.Range("G1:G37").ClearContents()

# --- Match the saved selection/active cell recorded in the workbook. ---
[void]$ws.Range("M45").Select()
